$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "18 to 19 years" column (E),
# shifting the old E:K columns to G:M.
$ws.Range("E:F").Insert()

# New header labels for the two inserted columns.
$ws.Range("E1").Value = "General warehousing and storage - Under 16 years"
$ws.Range("F1").Value = "General warehousing and storage - 16 to 17 years"

# Fill the new columns with a text "0" (not numeric 0) for every data row,
# matching the corrected VAR labeling in the source data.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = "'0"
    $ws.Cells.Item($r, 6).Value = "'0"
}
